# Generate Report for Handoff
#
# Updates the localization-status report to reflect a fresh handoff run:
#   - Status cells flip from "Handed back: in sync with en-US" to
#     "Ready for handoff"
#   - The associated timestamp cells are bumped to the new generation time
#   - The now-shorter Status column narrows accordingly (as Excel does when
#     the column is re-fit to its new contents)

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("E2").Value = "Ready for handoff"
$ws.Range("F2").Value = "Ready for handoff"
$ws.Range("G2").Value = "2016-08-22 15:09:24"
$ws.Columns("E:F").ColumnWidth = 16.33

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-22 15:09:18"
$ws.Columns("C:C").ColumnWidth = 16.33

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("H2").Value = "2016-08-22 15:09:24"
$ws.Columns("C:C").ColumnWidth = 16.33
